# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
# Copy the formatting of the existing last header cell (AC1, style s="1")
# onto the new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-46): Wins=66, Losses=96, Ties=0 for every row ---
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
